# B6-PowerPoint.pptx edit — Fri, Mar 20, 2020 11:06:42 AM
#
# 1) Re-style the three tables (slides 14, 15, 16) from the custom
#    "Table_0" style {CCE7A942-03B7-43F1-BE50-0050E77A4ED2} to the
#    built-in "No Style, No Grid" style {5E579BFD-AA1A-498D-B762-5F5B03FDCECE}.
# 2) Swap the deck's two themes: the slide master's theme ("Integral" /
#    Red Violet clrScheme) and the notes master's theme ("Office Theme")
#    trade color schemes, so the slide master ends up on the Office
#    Theme palette.

$p = $ppt.ActivePresentation

$oldTableStyle = "{CCE7A942-03B7-43F1-BE50-0050E77A4ED2}"
$newTableStyle = "{5E579BFD-AA1A-498D-B762-5F5B03FDCECE}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $table = $shape.Table
            # Every table in this deck currently carries the custom
            # "Table_0" style ($oldTableStyle); re-point it at the
            # built-in "No Style, No Grid" style.
            if ($table.Style -eq $oldTableStyle) {
                $table.ApplyStyle($newTableStyle)
            }
        }
    }
}

# --- Theme swap: clrScheme of theme1.xml (Integral / Red Violet) becomes
# the Office Theme palette that currently lives in theme2.xml. ---

function Set-ThemeColor($scheme, [int]$index, [string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    $scheme.Item($index).RGB = ($b * 65536) + ($g * 256) + $r
}

$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$masterColorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    Set-ThemeColor $masterColorScheme $i $officeColors[$i - 1]
}
